$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 9473.52
$ws.Range("B12").Value = 9390.8799999999992
$ws.Range("C12").Value = 105.78
$ws.Range("D12").Value = 106.71
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = 0.88
$ws.Range("G12").Value = 42620.766250000001
$ws.Range("G12").NumberFormat = "m/d/yy h:mm"
$ws.Range("H12").Value = $true
